$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 369  # was 370
$ws.Range("F5").Value = 1792  # was 1790
$ws.Range("F6").Value = 1792  # was 1790
$ws.Range("F7").Value = 394  # was 392
$ws.Range("F8").Value = 1469  # was 1467
$ws.Range("F9").Value = 839  # was 838
$ws.Range("F10").Value = 369  # was 367
$ws.Range("F11").Value = 718  # was 716
$ws.Range("F12").Value = 13059  # was 13047
$ws.Range("F13").Value = 12970  # was 12964
$ws.Range("F15").Value = 757  # was 755
$ws.Range("F17").Value = 539  # was 538
$ws.Range("F19").Value = 616  # was 615
$ws.Range("F20").Value = 2039  # was 2037
$ws.Range("F21").Value = 47  # was 46
$ws.Range("F22").Value = 23  # was 21
$ws.Range("F23").Value = 26  # was 25
$ws.Range("F25").Value = 157  # was 155
$ws.Range("F26").Value = 264  # was 263
$ws.Range("F27").Value = 721  # was 719

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 29  # was 28
$ws.Range("F4").Value = 14  # was 13
$ws.Range("F6").Value = 57  # was 56
$ws.Range("F7").Value = 93  # was 92
$ws.Range("F8").Value = 8  # was 5

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 369  # was 370
$ws.Range("F7").Value = 1792  # was 1790
$ws.Range("F8").Value = 1792  # was 1790
$ws.Range("F9").Value = 394  # was 392
$ws.Range("F10").Value = 1469  # was 1467
$ws.Range("F11").Value = 839  # was 838
$ws.Range("F12").Value = 369  # was 367
$ws.Range("F14").Value = 718  # was 716
$ws.Range("F15").Value = 13059  # was 13047
$ws.Range("F16").Value = 12970  # was 12964
$ws.Range("F18").Value = 757  # was 755
$ws.Range("F20").Value = 539  # was 538
$ws.Range("F22").Value = 616  # was 615
$ws.Range("F23").Value = 29  # was 28
$ws.Range("F24").Value = 14  # was 13
$ws.Range("F25").Value = 2039  # was 2037
$ws.Range("F26").Value = 47  # was 46
$ws.Range("F27").Value = 23  # was 21
$ws.Range("F28").Value = 26  # was 25
$ws.Range("F31").Value = 57  # was 56
$ws.Range("F32").Value = 157  # was 155
$ws.Range("F33").Value = 264  # was 263
$ws.Range("F34").Value = 721  # was 719
$ws.Range("F35").Value = 93  # was 92
$ws.Range("F36").Value = 8  # was 5
